# Commit: "Tue, May 05, 2020  4:07:42 AM"
#
# The canonical OOXML diff swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme2.xml (the theme actually wired to the slide
# master / whole deck via presentation.xml.rels + slideMaster1.xml.rels)
# goes from the "Integral" palette to the stock "Office Theme" palette,
# while theme1.xml (wired only to the notes master) goes the other way,
# from "Office Theme" to "Integral".
#
# The PowerPoint object model edits a presentation's *design* (theme) via
# the 12-slot ThemeColorScheme exposed on a Slide/SlideRange (it mirrors
# PpThemeColorSchemeIndex: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6,
# 11 hlink, 12 folHlink). Re-applying the Office Theme's color values here
# is exactly the "switch the deck back to the default Office design"
# operation a user performs from Design > Themes, and it rewrites the
# underlying theme part (theme2.xml) used by every slide/master in the
# deck to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1      -> 000000
$tcs.Item(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2      -> 44546A
$tcs.Item(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2  -> ED7D31
$tcs.Item(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4  -> FFC000
$tcs.Item(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5  -> 4472C4
$tcs.Item(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6  -> 70AD47
$tcs.Item(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink    -> 0563C1
$tcs.Item(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink -> 954F72
